$wb = $excel.ActiveWorkbook

# Report regenerated for archive: the file
# "13ed6402-d7e1-4f6e-9348-9283b488acdf.md" has moved from
# "Ready for handoff" back to "In Translation" status.

# --- Overview sheet: zh-cn (col B) and de-de (col C) status columns, row 5 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B5").Value = "In Translation"
$wsOverview.Range("C5").Value = "In Translation"

# --- zh-cn sheet: Status column (C), row 5 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C5").Value = "In Translation"

# --- de-de sheet: Status column (C), row 5 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C5").Value = "In Translation"
